$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2111
$ws.Range("F5").Value = 8019
$ws.Range("F7").Value = 7248
$ws.Range("G8").Value = 70
$ws.Range("F10").Value = 81
$ws.Range("F11").Value = 1202
$ws.Range("F13").Value = 210
$ws.Range("F14").Value = 582
$ws.Range("F16").Value = 48
$ws.Range("F19").Value = 41
$ws.Range("F20").Value = 1300
$ws.Range("F21").Value = 1280
$ws.Range("F24").Value = 1291
$ws.Range("F26").Value = 173
$ws.Range("F28").Value = 32
$ws.Range("F34").Value = 169
$ws.Range("F37").Value = 565
$ws.Range("F38").Value = 583
$ws.Range("F40").Value = 100
$ws.Range("F42").Value = 124
$ws.Range("F45").Value = 626

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 13
$ws.Range("F18").Value = 214
$ws.Range("F27").Value = 46
$ws.Range("F31").Value = 9
$ws.Range("F32").Value = 1013
$ws.Range("F36").Value = 132
$ws.Range("F39").Value = 112
$ws.Range("F40").Value = 151
$ws.Range("F41").Value = 24
$ws.Range("F42").Value = 18
$ws.Range("F47").Value = 3

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 712
$ws.Range("F7").Value = 223
$ws.Range("F8").Value = 114
$ws.Range("F9").Value = 1834
$ws.Range("F10").Value = 2734

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 712
$ws.Range("F6").Value = 8019
$ws.Range("F7").Value = 223
$ws.Range("F10").Value = 114
$ws.Range("F11").Value = 1834
$ws.Range("F12").Value = 2734
$ws.Range("F14").Value = 214
$ws.Range("F15").Value = 1202
$ws.Range("F18").Value = 582
$ws.Range("F20").Value = 1300
$ws.Range("F22").Value = 1280
$ws.Range("F24").Value = 1291
$ws.Range("F25").Value = 173
$ws.Range("F32").Value = 46
$ws.Range("F34").Value = 169
$ws.Range("F35").Value = 9
$ws.Range("F38").Value = 583
$ws.Range("F40").Value = 100
$ws.Range("F42").Value = 124
$ws.Range("F43").Value = 626
$ws.Range("F44").Value = 112
$ws.Range("F45").Value = 151

